# Update status for d7e5b1d5-4b03-43de-ad0e-29c77e9b1489.md from
# "Ready for handoff" back to "In Translation" on every report sheet
# (Overview + each locale sheet), as part of generating the archive report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B6").Value = "In Translation"
$overview.Range("C6").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B6").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B6").Value = "In Translation"
